$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 13; this pushes the existing rows 13-36
# down to 14-37, which is exactly the shift seen across the rest of the
# diff (every subsequent row's values equal what used to be one row above).
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the added weekly price record.
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "Vega Modelo de Temuco"
$ws.Range("C13").Value = "La Araucanía"
$ws.Range("D13").Value = 44690
$ws.Range("E13").Value = 9
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100108
$ws.Range("H13").Value = "Tropicales y subtropicales"
$ws.Range("I13").Value = 100108003
$ws.Range("J13").Value = "Maracuyá"
$ws.Range("K13").Value = "Sin especificar"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 25
$ws.Range("N13").Value = 34000
$ws.Range("O13").Value = 34000
$ws.Range("P13").Value = 34000
$ws.Range("Q13").Value = "$/caja 18 kilos"
$ws.Range("R13").Value = "Región de Arica y Parinacota"
$ws.Range("S13").Value = 1889
$ws.Range("T13").Value = 18
